# "break out stock.yaml completed"
# - Convert bsecode (column D) for rows 90-99 on the "day" sheet from text to numeric.
# - Append 8 new rows (100-107) of stock data to the "day" sheet, keeping bsecode (D)
#   as text for the newly appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Convert D90:D99 (bsecode) from text to real numbers ---
$bsecodes = @{
    90 = 500290
    91 = 500488
    92 = 500825
    93 = 539524
    94 = 533398
    95 = 532424
    96 = 500800
    97 = 540611
    98 = 535755
    99 = 532720
}

foreach ($r in $bsecodes.Keys) {
    $ws.Cells.Item($r, 4).Value = $bsecodes[$r]
}

# --- Append new rows 100-107 ---
# Columns: A=sr, B=nsecode, C=name, D=bsecode, E=per_chg, F=close, G=volume, H=timeframe, I=Date Time
$newRows = @(
    @{ Row=100; Sr=1; Nse="BAJFINANCE"; Name="Bajaj Finance Limited";                       Bse="500034"; Chg=-0.43; Close=7068.05;  Vol=973503;   Dt="09/07/2024 11:35:34" },
    @{ Row=101; Sr=2; Nse="POLYCAB";    Name="Polycab India Ltd";                           Bse="542652"; Chg=-1.41; Close=6454.15;  Vol=547696;   Dt="09/07/2024 11:35:34" },
    @{ Row=102; Sr=3; Nse="CUMMINSIND"; Name="Cummins India Limited";                       Bse="500480"; Chg=-0.38; Close=3985.8;   Vol=432053;   Dt="09/07/2024 11:35:34" },
    @{ Row=103; Sr=4; Nse="HAVELLS";    Name="Havells India Limited";                       Bse="517354"; Chg=1.81;  Close=1921.05;  Vol=1106258;  Dt="09/07/2024 11:35:34" },
    @{ Row=104; Sr=5; Nse="EXIDEIND";   Name="Exide Industries Limited";                    Bse="500086"; Chg=1.23;  Close=578.45;   Vol=7274064;  Dt="09/07/2024 11:35:34" },
    @{ Row=105; Sr=6; Nse="JUBLFOOD";   Name="Jubilant Foodworks Limited";                  Bse="533155"; Chg=-0.13; Close=573.4;    Vol=1888512;  Dt="09/07/2024 11:35:34" },
    @{ Row=106; Sr=7; Nse="BPCL";       Name="Bharat Petroleum Corporation Limited";        Bse="500547"; Chg=0.23;  Close=300.2;    Vol=18094862; Dt="09/07/2024 11:35:34" },
    @{ Row=107; Sr=8; Nse="ASHOKLEY";   Name="Ashok Leyland Limited";                       Bse="500477"; Chg=0.99;  Close=228.28;   Vol=13248300; Dt="09/07/2024 11:35:34" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.Sr
    $ws.Cells.Item($r, 2).Value = $row.Nse
    $ws.Cells.Item($r, 3).Value = $row.Name

    # bsecode stays as text in the new rows, so force text formatting, assign,
    # then strip the formatting again so no style index is attached to the cell.
    $cellD = $ws.Cells.Item($r, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $row.Bse
    $cellD.ClearFormats()

    $ws.Cells.Item($r, 5).Value = $row.Chg
    $ws.Cells.Item($r, 6).Value = $row.Close
    $ws.Cells.Item($r, 7).Value = $row.Vol
    $ws.Cells.Item($r, 8).Value = "day"
    $ws.Cells.Item($r, 9).Value = $row.Dt
}
